$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17, shifting existing rows 17-49 down to 18-50.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with a new weekly price record (same market/category
# as the surrounding rows, new date and price figures).
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "Terminal La Palmera de La Serena"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44498
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 100112052
$ws.Range("G17").Value = "Albahaca"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 900
$ws.Range("K17").Value = 3800
$ws.Range("L17").Value = 4000
$ws.Range("M17").Value = 3900
$ws.Range("N17").Value = "$/paquete"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 3900
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"
